$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.150.99'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '3.430.21'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.71'
$ws.Range('E5').Value = '  -1.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.03'
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.579'
$ws.Range('E8').Value = '  -2.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.127'
$ws.Range('E9').Value = '  -6.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.05'
$ws.Range('E10').Value = '  -3.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.420'
$ws.Range('E11').Value = '  -3.75%  '
$ws.Range('D12').Value = '4.044.81'
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.134'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.20'
$ws.Range('E14').Value = '  +4.89%  '
$ws.Range('D15').Value = '66.134.87'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000173'
$ws.Range('E16').Value = '  -5.31%  '
$ws.Range('D17').Value = '3.450.68'
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.16'
$ws.Range('E18').Value = '  -3.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.16'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '383.70'
$ws.Range('E20').Value = '  -3.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.75'
$ws.Range('E21').Value = '  -3.16%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.27'
$ws.Range('E23').Value = '  -1.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.69'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.527'
$ws.Range('E25').Value = '  -2.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000118'
$ws.Range('E26').Value = '  -4.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.04'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.174'
$ws.Range('E28').Value = '  -4.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.00'
$ws.Range('E30').Value = '  -5.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.39'
$ws.Range('E31').Value = '  -5.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.01'
$ws.Range('E32').Value = '  -3.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.13'
$ws.Range('E33').Value = '  -4.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.17'
$ws.Range('E34').Value = '  -3.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.22'
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.856'
$ws.Range('E37').Value = '  -4.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.88'
$ws.Range('E38').Value = '  -2.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.79'
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.55'
$ws.Range('E40').Value = '  -4.09%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.51'
$ws.Range('E41').Value = '  -2.88%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.759.01'
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '25.60'
$ws.Range('E43').Value = '  -3.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0710'
$ws.Range('E44').Value = '  -5.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.09'
$ws.Range('E45').Value = '  -2.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.47'
$ws.Range('E46').Value = '  -6.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0294'
$ws.Range('E47').Value = '  -5.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '332.24'
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.04'
$ws.Range('E49').Value = '  -5.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '32.82'
$ws.Range('E50').Value = '  -2.78%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.102'
$ws.Range('E51').Value = '  -2.72%  '
